# Update gh-pages output (杭州-漫展信息.xlsx) to the regenerated scrape.
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both carry the
# same underlying rows, so most edits are mirrored between them. One row
# ("杭州·幻想物语动漫游戏展") was cancelled: on "展览" it stays in place
# (title suffixed, marked 不可售) while on "全部类型" it is dropped and the
# rows below it shift up by one, with a brand-new row ("杭州·第九届萌次元
# 动漫嘉年华") appended after the shifted block to keep row 18+ untouched.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $text)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

function Set-NumCell {
    param($ws, $row, $col, $num)
    $ws.Cells.Item($row, $col).Value = $num
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Plain "想去人数" (F column) refreshes - title/location/etc unchanged.
$s1_f = @{
    2 = 10160; 4 = 2524; 6 = 285; 9 = 761; 13 = 3158; 14 = 2363; 16 = 2090; 17 = 2090;
    22 = 558; 24 = 240; 25 = 7; 26 = 16; 29 = 370; 31 = 363; 32 = 581; 33 = 49; 34 = 234;
    35 = 2; 36 = 1571; 37 = 49; 38 = 330; 39 = 1670; 40 = 113; 41 = 426; 42 = 49; 43 = 439; 44 = 958
}
foreach ($row in $s1_f.Keys) {
    Set-NumCell $ws1 $row 6 $s1_f[$row]
}

# Row 11: 杭州·幻想物语动漫游戏展 -> cancelled, no longer sellable.
Set-TextCell $ws1 11 3 "杭州·幻想物语动漫游戏展（取消）"
Set-NumCell  $ws1 11 6 1226
Set-TextCell $ws1 11 7 "不可售"

# ---------------------------------------------------------------------
# Sheet 2: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# Rows untouched by the row-12 cancellation (same row numbers as before) -
# plain F-column refreshes.
$s4_f = @{
    2 = 10160; 4 = 2524; 8 = 285; 11 = 761;
    21 = 558; 23 = 240; 24 = 7; 25 = 16; 28 = 370; 30 = 363; 31 = 581; 35 = 49; 36 = 234;
    37 = 2; 38 = 1571; 39 = 49; 41 = 330; 42 = 1670; 43 = 113; 45 = 426; 46 = 49; 47 = 439; 48 = 958
}
foreach ($row in $s4_f.Keys) {
    Set-NumCell $ws4 $row 6 $s4_f[$row]
}

# Row 12 (杭州·幻想物语动漫游戏展) is removed from this sheet: rows 13-17
# shift up into 12-16 (each also picking up its own refreshed numbers),
# and a brand-new row is written into the vacated slot 17 so row 18
# (杭州·ESCC电竞博览会 倒霉死勒内场票) keeps its original row number.
# NOTE: column A (the 0-based sequence number) is intentionally left
# untouched on every one of these rows.

# New row 12 <- old row 13 content (杭州·排球少年only·春日校庆)
Set-TextCell $ws4 12 2 "2024-03-30"
Set-TextCell $ws4 12 3 "杭州·排球少年only·春日校庆"
Set-TextCell $ws4 12 4 "之江路149号 云栖培训基地"
Set-TextCell $ws4 12 5 "2024.03.30 10:00-03.31 17:00"
Set-NumCell  $ws4 12 6 1050
Set-NumCell  $ws4 12 7 89
Set-TextCell $ws4 12 8 "https://show.bilibili.com/platform/detail.html?id=81511"
Set-TextCell $ws4 12 9 "//i0.hdslb.com/bfs/openplatform/202402/RDI807mS1708410823039.jpeg"

# New row 13 <- old row 14 content (杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会)
Set-TextCell $ws4 13 2 "2024-04-04"
Set-TextCell $ws4 13 3 "杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会"
Set-TextCell $ws4 13 4 "钱江世纪城奔竞大道353号 杭州国际博览中心"
Set-TextCell $ws4 13 5 "2024.04.04 09:30-04.05 16:30"
Set-NumCell  $ws4 13 6 3158
Set-NumCell  $ws4 13 7 75
Set-TextCell $ws4 13 8 "https://show.bilibili.com/platform/detail.html?id=81450"
Set-TextCell $ws4 13 9 "//i1.hdslb.com/bfs/openplatform/202403/OfpkJ50P1709548942017.png"

# New row 14 <- old row 15 content (杭州·ELECTRIC COMIC动漫游戏展)
Set-TextCell $ws4 14 2 "2024-04-04"
Set-TextCell $ws4 14 3 "杭州·ELECTRIC COMIC动漫游戏展"
Set-TextCell $ws4 14 4 "望江东路333号 杭州瑞莱克斯大酒店"
Set-TextCell $ws4 14 5 "2024.04.04 10:00-04.05 17:00"
Set-NumCell  $ws4 14 6 2363
Set-NumCell  $ws4 14 7 63
Set-TextCell $ws4 14 8 "https://show.bilibili.com/platform/detail.html?id=82270"
Set-TextCell $ws4 14 9 "//i2.hdslb.com/bfs/openplatform/202403/JmFXyFgc1710844373405.jpeg"

# New row 15 <- old row 16 content (杭州·梦漫星河动漫展, first occurrence)
Set-TextCell $ws4 15 2 "2024-04-04"
Set-TextCell $ws4 15 3 "杭州·梦漫星河动漫展"
Set-TextCell $ws4 15 4 "德胜东路2539号 梦马汽车小镇"
Set-TextCell $ws4 15 5 "2024.04.04 10:00-04.05 17:00"
Set-NumCell  $ws4 15 6 2090
Set-NumCell  $ws4 15 7 58.5
Set-TextCell $ws4 15 8 "https://show.bilibili.com/platform/detail.html?id=81699"
Set-TextCell $ws4 15 9 "//i0.hdslb.com/bfs/openplatform/202402/sZfZd47Y1706868453434.jpeg"

# New row 16 <- old row 17 content (杭州·梦漫星河动漫展, second occurrence);
# only the 想去人数 count changes here.
Set-NumCell $ws4 16 6 2090

# New row 17: brand-new entry, 杭州·第九届萌次元动漫嘉年华 (not sellable).
Set-TextCell $ws4 17 2 "2024-04-04"
Set-TextCell $ws4 17 3 "杭州·第九届萌次元动漫嘉年华"
Set-TextCell $ws4 17 4 "长乐路29号五组2幢 杭州运河文化发布中心"
Set-TextCell $ws4 17 5 "2024.04.04 10:00-04.05 17:00"
Set-NumCell  $ws4 17 6 246
Set-TextCell $ws4 17 7 "不可售"
Set-TextCell $ws4 17 8 "https://show.bilibili.com/platform/detail.html?id=78866"
Set-TextCell $ws4 17 9 "//i1.hdslb.com/bfs/openplatform/202311/8jSeAOZH1700636327971.jpeg"

Write-Output "edits applied"
